# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 42 (pushing existing rows 42-107 down to 43-108)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 42; Excel shifts rows 42:107 down to 43:108.
$ws.Rows("42:42").Insert()

# Populate every column of the new row 42 with its values.
$ws.Range("A42").Value = 3
$ws.Range("B42").Value = "Femacal de La Calera"
$ws.Range("C42").Value = "Coquimbo"
$ws.Range("D42").Value = 44536
$ws.Range("D42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E42").Value = 5
$ws.Range("F42").Value = 100112026
$ws.Range("G42").Value = "Haba"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 125
$ws.Range("K42").Value = 8000
$ws.Range("L42").Value = 8000
$ws.Range("M42").Value = 8000
$ws.Range("N42").Value = "$/saco 25 kilos"
$ws.Range("O42").Value = "Provincia de Limarí"
$ws.Range("P42").Value = 320
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = "Hortaliza"
